$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the trailing paragraphs that the commit deletes entirely:
#    "git add -all" / (blank) / "Zatim je napravljen..." / (blank) /
#    "git add -all" / "git commit"
#    These are Paragraphs 3..8 in the original document (1-based).
# ---------------------------------------------------------------------------
if ($d.Paragraphs.Count -ge 8) {
    $startPara = $d.Paragraphs.Item(3)
    $endPara   = $d.Paragraphs.Item(8)
    $killRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $killRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Merge paragraph 1 (the long sentence) with paragraph 2 (the bookmark-only
#    paragraph) by deleting the paragraph mark that separates them, so the
#    bookmark ends up inside the same paragraph as the text, as in the diff.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$mark = $d.Range($p1.Range.End - 1, $p1.Range.End)
$mark.Delete()

# ---------------------------------------------------------------------------
# 3) Rewrite the sentence text, in left-to-right order. Each Find/Replace
#    below targets either a whole run or a span that starts/ends exactly on
#    a spellStart/spellEnd proofing-mark boundary, so the proofErr markers
#    stay balanced (no orphan tags left behind).
# ---------------------------------------------------------------------------
$rng = $d.Content

# "Inicijalizacija" -> "Prvi commit, za"
$rng.Find.Execute("Inicijalizacija", $true, $false, $false, $false, $false, $true, 1, $false, "Prvi commit, za", 2) | Out-Null

# "brancha" -> "lekcije"
$rng.Find.Execute("brancha", $true, $false, $false, $false, $false, $true, 1, $false, "lekcije", 2) | Out-Null

# "LearnigMaterial" -> "koje"
$rng.Find.Execute("LearnigMaterial", $true, $false, $false, $false, $false, $true, 1, $false, "koje", 2) | Out-Null

# ", u " -> " se "
$rng.Find.Execute(", u ", $true, $false, $false, $false, $false, $true, 1, $false, " se ", 2) | Out-Null

# "folderu" -> "ticu"
$rng.Find.Execute("folderu", $true, $false, $false, $false, $false, $true, 1, $false, "ticu", 2) | Out-Null

# " je " -> " "  (first occurrence only: search forward from here, no wrap)
$rng.Find.Execute(" je ", $true, $false, $false, $false, $false, $true, 0, $false, " ", 2) | Out-Null

# "obrisan" -> "programiranja"
$rng.Find.Execute("obrisan", $true, $false, $false, $false, $false, $true, 0, $false, "programiranja", 2) | Out-Null

# "stari ReadMe fajl, promena je dodata" -> "WPpluginove"
$rng.Find.Execute("stari ReadMe fajl, promena je dodata", $true, $false, $false, $false, $false, $true, 0, $false, "WPpluginove", 2) | Out-Null

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host ("Para " + $i + ": [" + $p.Range.Text + "]")
}
